# "add complete report agg and detail"
# Adds a new "order" field to the categories table schema (leader_direct
# sheet, which documents table/field metadata) and a corresponding
# "order" data column to the categories sheet itself, renumbering the
# existing sort-order values along the way.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) leader_direct: insert a new schema row describing the "order"
#    field of the categories table, just above the existing "status"
#    field row (old row 135 -> becomes row 136).
# ------------------------------------------------------------------
$lead = $wb.Worksheets.Item("leader_direct")

# Insert a blank row at 135, pushing the old row 135 ("status") down
# to row 136, then copy the (soon to be) row-136 formatting onto the
# freshly inserted row 135 so the new row matches its neighbours.
$lead.Range("A135:P135").Insert(-4121)
$lead.Range("A136:P136").Copy()
$lead.Range("A135:P135").PasteSpecial(-4122)

$lead.Range("A135").Value = "categories"
$lead.Range("B135").Value = "order"
$lead.Range("D135").Value = "Thứ tự sắp xếp"
$lead.Range("C135").Value = "Thứ Tự"
$lead.Range("E135").Value = "INTEGER"
$lead.Range("M135").Value = 1

# Renumber the order_1 (P) column so it is sequential 1..135 for every
# data row (some rows had drifted out of sequence before this edit,
# and the newly inserted/shifted rows need values too).
for ($i = 75; $i -le 136; $i++) {
    $lead.Range("P$i").Value = $i - 1
}

# ------------------------------------------------------------------
# 2) categories: insert a new "order" column before the existing
#    "status" column (M -> N), and fill it with sequential values.
# ------------------------------------------------------------------
$cat = $wb.Worksheets.Item("categories")

$cat.Range("M1:M55").Insert(-4161)
$cat.Columns.Item(13).ColumnWidth = 12.67

$cat.Range("M1").Value = "order"
for ($i = 2; $i -le 55; $i++) {
    $cat.Range("M$i").Value = $i - 1
}

# ------------------------------------------------------------------
# 3) View state: leader_direct becomes the active/selected sheet
#    (instead of categories), with the cursor on the newly added row.
# ------------------------------------------------------------------
$cat.Range("K7").Select()

$lead.Activate()
$lead.Range("C136").Select()
